$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.054.97'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '1.893.18'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.53'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5217'
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3760'
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07257'
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.09'
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8983'
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08181'
$ws.Range("E12").Value = '  +3.45%  '
$ws.Range("D13").Value = '1.938.43'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.32'
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.296'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.56'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").Value = '27.083.77'
$ws.Range("E20").Value = '  -0.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.076'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.405'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.75'
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.287'
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.18'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.733'
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.07'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.785'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.845'
$ws.Range("E30").Value = '  -2.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09216'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05028'
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7874'
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.208'
$ws.Range("E34").Value = '  -3.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.424'
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.970'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.604'
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5704'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.075'
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.016'
$ws.Range("E41").Value = '  +0.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.554'
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.14'
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1512'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4854'
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.02'
$ws.Range("E47").Value = '  -1.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.620'
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.14'
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("E51").Value = '  -0.16%  '
